# Generate Report for Handback
# The localization file "042c2ff2-4e04-4c4c-92af-fd9a54fc9367.md" has been
# handed back and is now in sync with en-US. Update the zh-cn and de-de
# status sheets with the new handback file / datetime, and (for de-de)
# flip the Status column to reflect the handback.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$targetFile = "042c2ff2-4e04-4c4c-92af-fd9a54fc9367.md"
$targetUrl  = "https://github.com/OpenLocalizationTestOrg/oltest/blob/e2b6fe8f014711238b888fe662cc7baea54d5799/e2e/042c2ff2-4e04-4c4c-92af-fd9a54fc9367.md"

# ---------------------------------------------------------------------
# zh-cn sheet, row 2 (042c2ff2-....md)
# ---------------------------------------------------------------------
$zhcn.Range("I2").Value = $targetFile
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $targetUrl, [System.Type]::Missing, [System.Type]::Missing, $targetFile) | Out-Null

$zhcn.Range("J2").Value = "042c2ff2-4e04-4c4c-92af-fd9a54fc9367.c5733616c1a95f3f23cca384c1e46e0aee842dd9.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-14 16:53:58"

# ---------------------------------------------------------------------
# de-de sheet, row 2 (042c2ff2-....md)
# ---------------------------------------------------------------------
$dede.Range("C2").Value = "Handed back: in sync with en-US"

$dede.Range("I2").Value = $targetFile
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276
$dede.Hyperlinks.Add($dede.Range("I2"), $targetUrl, [System.Type]::Missing, [System.Type]::Missing, $targetFile) | Out-Null

$dede.Range("J2").Value = "042c2ff2-4e04-4c4c-92af-fd9a54fc9367.c5733616c1a95f3f23cca384c1e46e0aee842dd9.de-de.xlf"
$dede.Range("K2").Value = "2016-08-14 16:54:11"

# ---------------------------------------------------------------------
# Widen columns that now hold longer content (matches Excel's autofit
# behaviour after the handback data was filled in).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 39.4574236188616
$zhcn.Columns.Item(10).ColumnWidth = 40

$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 39.4574236188616
$dede.Columns.Item(10).ColumnWidth = 40
